$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- E2: saldo awal tweak ---
$ws.Range("E2").Value = 685525

# --- D3: add an extra expense to the existing formula ---
$ws.Range("D3").Formula = "=60000+260000"

# --- D4: append more entries to the existing formula ---
$ws.Range("D4").Formula = "=2877500+1537000+3649500+5980000+2600000+41700000+498000"

# New shared strings must be interned in the same order the original
# workbook introduced them (60..65), so set the B-column labels first in
# that order before touching the numeric/formula cells.
$ws.Range("B5").Value = "BELI kresek"
$ws.Range("B6").Value = "A/R"
$ws.Range("B8").Value = "SALES - cash/retail"
$ws.Range("B7").Value = "BENSIN - RUSH"
$ws.Range("B9").Value = "SELISIH - lebih"
$ws.Range("B10").Value = "SETOR KE BANK"

# --- Row 5: BELI kresek ---
$ws.Range("D5").Value = 100000

# --- Row 6: A/R ---
$ws.Range("C6").Formula = "=25000000+16700000+25577000"

# --- Row 7: BENSIN - RUSH ---
$ws.Range("D7").Value = 200000

# --- Row 8: SALES - cash/retail ---
$ws.Range("C8").Formula = "=16816475+22816525-25577000"

# --- Row 9: SELISIH - lebih ---
$ws.Range("C9").Value = 2000

# --- Row 10: SETOR KE BANK ---
$ws.Range("D10").Formula = "=22000000"

# --- Row 11: new day, Wages Expense ---
$ws.Range("A11").Value = 44313
$ws.Range("B11").Value = "Wages Expense"
$ws.Range("D11").Value = 45000

# --- Row 12: A/R ---
$ws.Range("B12").Value = "A/R"
$ws.Range("C12").Formula = "=14625000+360000+28500000+8800000"

# --- Row 13: TRANSFER BCA ---
$ws.Range("B13").Value = "TRANSFER BCA"
$ws.Range("D13").Formula = "=14985000+600000+569000+175000+8800000"

# --- Move the active selection to D14 ---
$ws.Range("D14").Select()
